$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")

# Add the new row of field test data (row 8)
$ws.Range("A8").Value = 45793
$ws.Range("B8").Value = "PRESENCE"
$ws.Range("C8").Value = 0.30208333333333331
$ws.Range("D8").Value = 0.58680555555555558
$ws.Range("E8").Value = 18
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = "Sunny, mild"
$ws.Range("H8").Value = $true
$ws.Range("I8").Value = "5 minutes 30 seconds"
$ws.Range("J8").Value = 330
$ws.Range("K8").Value = "Primary sweeps"
$ws.Range("L8").Value = "Worked downhill. Koda got onto odour 1/4 way through and pulled immediately downhill to about 1/2 along the transect. Took me 10 seconds to catch up!"

$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("C8:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update sheet view (scroll/selection state)
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("J9").Select()
